$d = $word.ActiveDocument

# Remove the horizontal-rule paragraphs (w:pict / v:rect "hr" placeholders).
# These paragraphs contain no real text -- just the paragraph mark -- so we
# can identify them reliably by Range.Text being a lone CR (chr 13) and
# delete each one's Range, which removes the paragraph (including its
# paragraph mark) entirely from the document body.
# Walk backwards so deleting a paragraph never invalidates the index of a
# paragraph we haven't visited yet.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt.Length -eq 1 -and [int]$txt[0] -eq 13) {
        $p.Range.Delete()
    }
}
